# Edit slide 4 ("Challenges" slide): merge the trailing " API/AT mode" and
# ") " runs in the "Resolving unexpected behavior..." bullet into a single
# run with text " API/AT mode) ", matching the author's final touch-up.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# Find the paragraph that contains the "Xbee's ... API/AT mode)" bullet.
$target = $null
for ($i = 1; $i -le $tr.Paragraphs().Count; $i++) {
    $para = $tr.Paragraphs($i)
    if ($para.Text.IndexOf("API/AT mode") -ge 0) {
        $target = $para
        break
    }
}

$ptext = $target.Text
$oldChunk = " API/AT mode) "
$localIdx = $ptext.IndexOf($oldChunk)
$startPos = $target.Start + $localIdx
$len = $oldChunk.Length

$chars = $tr.Characters($startPos, $len)
$chars.Text = " API/AT mode) "
